$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 37.84230566666667
$ws.Range("H2").Value = 113.526917
$ws.Range("I2").Value = 0.5048163536019187
$ws.Range("J2").Value = 0.5048163536019187
$ws.Range("M2").Value = 20.94338233333333
$ws.Range("N2").Value = 62.830147
$ws.Range("O2").Value = 0.06014699761632732
$ws.Range("P2").Value = 0.06014699761632732
$ws.Range("Q2").Value = 792.5458759518665
$ws.Range("R2").Value = 7132.912883566799
$ws.Range("S2").Value = 0.03036318801677765
$ws.Range("T2").Value = 0.03036318801677765

# Row 3
$ws.Range("G3").Value = 37.84230566666667
$ws.Range("H3").Value = 113.526917
$ws.Range("I3").Value = 0.5048163536019187
$ws.Range("J3").Value = 0.5048163536019187
$ws.Range("M3").Value = 36.67092
$ws.Range("N3").Value = 110.01276
$ws.Range("O3").Value = 0.105314686172636
$ws.Range("P3").Value = 0.105314686172636
$ws.Range("Q3").Value = 1387.71216371788
$ws.Range("R3").Value = 12489.40947346092
$ws.Range("S3").Value = 0.05316457585440052
$ws.Range("T3").Value = 0.05316457585440053

# Row 4
$ws.Range("G4").Value = 37.84230566666667
$ws.Range("H4").Value = 113.526917
$ws.Range("I4").Value = 0.5048163536019187
$ws.Range("J4").Value = 0.5048163536019187
$ws.Range("M4").Value = 284.2332763333333
$ws.Range("N4").Value = 852.6998289999999
$ws.Range("O4").Value = 0.8162854462572833
$ws.Range("P4").Value = 0.8162854462572834
$ws.Range("Q4").Value = 10756.04252364413
$ws.Range("R4").Value = 96804.38271279717
$ws.Range("S4").Value = 0.4120742424779167
$ws.Range("T4").Value = 0.4120742424779167

# Row 5
$ws.Range("G5").Value = 37.84230566666667
$ws.Range("H5").Value = 113.526917
$ws.Range("I5").Value = 0.5048163536019187
$ws.Range("J5").Value = 0.5048163536019187
$ws.Range("M5").Value = 6.355709333333334
$ws.Range("N5").Value = 19.067128
$ws.Range("O5").Value = 0.01825286995375338
$ws.Range("P5").Value = 0.01825286995375338
$ws.Range("Q5").Value = 240.5146953204863
$ws.Range("R5").Value = 2164.632257884376
$ws.Range("S5").Value = 0.009214347252823803
$ws.Range("T5").Value = 0.009214347252823805

# Row 6
$ws.Range("G6").Value = 15.88630666666666
$ws.Range("H6").Value = 47.65891999999999
$ws.Range("I6").Value = 0.2119233292577262
$ws.Range("J6").Value = 0.2119233292577262
$ws.Range("M6").Value = 20.94338233333333
$ws.Range("N6").Value = 62.830147
$ws.Range("O6").Value = 0.06014699761632732
$ws.Range("P6").Value = 0.06014699761632732
$ws.Range("Q6").Value = 332.7129943845821
$ws.Range("R6").Value = 2994.41694946124
$ws.Range("S6").Value = 0.0127465519797086
$ws.Range("T6").Value = 0.01274655197970861

# Row 7
$ws.Range("G7").Value = 15.88630666666666
$ws.Range("H7").Value = 47.65891999999999
$ws.Range("I7").Value = 0.2119233292577262
$ws.Range("J7").Value = 0.2119233292577262
$ws.Range("M7").Value = 36.67092
$ws.Range("N7").Value = 110.01276
$ws.Range("O7").Value = 0.105314686172636
$ws.Range("P7").Value = 0.105314686172636
$ws.Range("Q7").Value = 582.5654808687998
$ws.Range("R7").Value = 5243.089327819199
$ws.Range("S7").Value = 0.02231863891343764
$ws.Range("T7").Value = 0.02231863891343765

# Row 8
$ws.Range("G8").Value = 15.88630666666666
$ws.Range("H8").Value = 47.65891999999999
$ws.Range("I8").Value = 0.2119233292577262
$ws.Range("J8").Value = 0.2119233292577262
$ws.Range("M8").Value = 284.2332763333333
$ws.Range("N8").Value = 852.6998289999999
$ws.Range("O8").Value = 0.8162854462572833
$ws.Range("P8").Value = 0.8162854462572834
$ws.Range("Q8").Value = 4515.416992702741
$ws.Range("R8").Value = 40638.75293432467
$ws.Range("S8").Value = 0.1729899293954722
$ws.Range("T8").Value = 0.1729899293954722

# Row 9
$ws.Range("G9").Value = 15.88630666666666
$ws.Range("H9").Value = 47.65891999999999
$ws.Range("I9").Value = 0.2119233292577262
$ws.Range("J9").Value = 0.2119233292577262
$ws.Range("M9").Value = 6.355709333333334
$ws.Range("N9").Value = 19.067128
$ws.Range("O9").Value = 0.01825286995375338
$ws.Range("P9").Value = 0.01825286995375338
$ws.Range("Q9").Value = 100.9687475535289
$ws.Range("R9").Value = 908.7187279817601
$ws.Range("S9").Value = 0.003868208969107734
$ws.Range("T9").Value = 0.003868208969107736

# Row 10
$ws.Range("G10").Value = 18.76675533333333
$ws.Range("H10").Value = 56.300266
$ws.Range("I10").Value = 0.2503485141672444
$ws.Range("J10").Value = 0.2503485141672445
$ws.Range("M10").Value = 20.94338233333333
$ws.Range("N10").Value = 62.830147
$ws.Range("O10").Value = 0.06014699761632732
$ws.Range("P10").Value = 0.06014699761632732
$ws.Range("Q10").Value = 393.0393321021224
$ws.Range("R10").Value = 3537.353988919102
$ws.Range("S10").Value = 0.01505771148486833
$ws.Range("T10").Value = 0.01505771148486834

# Row 11
$ws.Range("G11").Value = 18.76675533333333
$ws.Range("H11").Value = 56.300266
$ws.Range("I11").Value = 0.2503485141672444
$ws.Range("J11").Value = 0.2503485141672445
$ws.Range("M11").Value = 36.67092
$ws.Range("N11").Value = 110.01276
$ws.Range("O11").Value = 0.105314686172636
$ws.Range("P11").Value = 0.105314686172636
$ws.Range("Q11").Value = 688.1941834882399
$ws.Range("R11").Value = 6193.747651394159
$ws.Range("S11").Value = 0.02636537520330907
$ws.Range("T11").Value = 0.02636537520330908

# Row 12
$ws.Range("G12").Value = 18.76675533333333
$ws.Range("H12").Value = 56.300266
$ws.Range("I12").Value = 0.2503485141672444
$ws.Range("J12").Value = 0.2503485141672445
$ws.Range("M12").Value = 284.2332763333333
$ws.Range("N12").Value = 852.6998289999999
$ws.Range("O12").Value = 0.8162854462572833
$ws.Range("P12").Value = 0.8162854462572834
$ws.Range("Q12").Value = 5334.13635453939
$ws.Range("R12").Value = 48007.22719085451
$ws.Range("S12").Value = 0.2043558486068569
$ws.Range("T12").Value = 0.204355848606857

# Row 13
$ws.Range("G13").Value = 18.76675533333333
$ws.Range("H13").Value = 56.300266
$ws.Range("I13").Value = 0.2503485141672444
$ws.Range("J13").Value = 0.2503485141672445
$ws.Range("M13").Value = 6.355709333333334
$ws.Range("N13").Value = 19.067128
$ws.Range("O13").Value = 0.01825286995375338
$ws.Range("P13").Value = 0.01825286995375338
$ws.Range("Q13").Value = 119.2760420284498
$ws.Range("R13").Value = 1073.484378256048
$ws.Range("S13").Value = 0.004569578872210098
$ws.Range("T13").Value = 0.0045695788722101

# Row 14
$ws.Range("G14").Value = 2.467151666666667
$ws.Range("H14").Value = 7.401455
$ws.Range("I14").Value = 0.03291180297311068
$ws.Range("J14").Value = 0.03291180297311068
$ws.Range("M14").Value = 20.94338233333333
$ws.Range("N14").Value = 62.830147
$ws.Range("O14").Value = 0.06014699761632732
$ws.Range("P14").Value = 0.06014699761632732
$ws.Range("Q14").Value = 51.67050062932055
$ws.Range("R14").Value = 465.034505663885
$ws.Range("S14").Value = 0.001979546134972722
$ws.Range("T14").Value = 0.001979546134972723

# Row 15
$ws.Range("G15").Value = 2.467151666666667
$ws.Range("H15").Value = 7.401455
$ws.Range("I15").Value = 0.03291180297311068
$ws.Range("J15").Value = 0.03291180297311068
$ws.Range("M15").Value = 36.67092
$ws.Range("N15").Value = 110.01276
$ws.Range("O15").Value = 0.105314686172636
$ws.Range("P15").Value = 0.105314686172636
$ws.Range("Q15").Value = 90.47272139619999
$ws.Range("R15").Value = 814.2544925658
$ws.Range("S15").Value = 0.00346609620148878
$ws.Range("T15").Value = 0.003466096201488781

# Row 16
$ws.Range("G16").Value = 2.467151666666667
$ws.Range("H16").Value = 7.401455
$ws.Range("I16").Value = 0.03291180297311068
$ws.Range("J16").Value = 0.03291180297311068
$ws.Range("M16").Value = 284.2332763333333
$ws.Range("N16").Value = 852.6998289999999
$ws.Range("O16").Value = 0.8162854462572833
$ws.Range("P16").Value = 0.8162854462572834
$ws.Range("Q16").Value = 701.2466014279105
$ws.Range("R16").Value = 6311.219412851195
$ws.Range("S16").Value = 0.02686542577703743
$ws.Range("T16").Value = 0.02686542577703744

# Row 17
$ws.Range("G17").Value = 2.467151666666667
$ws.Range("H17").Value = 7.401455
$ws.Range("I17").Value = 0.03291180297311068
$ws.Range("J17").Value = 0.03291180297311068
$ws.Range("M17").Value = 6.355709333333334
$ws.Range("N17").Value = 19.067128
$ws.Range("O17").Value = 0.01825286995375338
$ws.Range("P17").Value = 0.01825286995375338
$ws.Range("Q17").Value = 15.68049887458222
$ws.Range("R17").Value = 141.12448987124
$ws.Range("S17").Value = 0.000600734859611743
$ws.Range("T17").Value = 0.0006007348596117432

